$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Inventory)
$ws.Range("C4").Value = 456000000.0
$ws.Range("D4").Value = 662000000.0
$ws.Range("E4").Value = 665000000.0
$ws.Range("F4").Value = 650000000.0
$ws.Range("G4").Value = 640000000.0

# Row 12 (Accounts Payable)
$ws.Range("C12").Value = 481000000.0
$ws.Range("D12").Value = 777000000.0
$ws.Range("E12").Value = 733000000.0
$ws.Range("F12").Value = 745000000.0
$ws.Range("G12").Value = 766000000.0

# Row 38 (Net Debt) - B38 changes from empty inline string to number
$ws.Range("B38").Value = 1267900000.0

# Row 39 (Total Debt) - B39 changes from empty inline string to number
$ws.Range("B39").Value = 2567500000.0
